$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric value changes
$ws.Range("B7").Value = 34.68740776195546
$ws.Range("B20").Value = 6558705.350000001
$ws.Range("B24").Value = 746482.5
$ws.Range("B27").Value = 99531000
$ws.Range("B33").Value = 116979388.0899999
$ws.Range("B35").Value = 4673469.8
$ws.Range("B36").Value = 720790.7700000001
$ws.Range("B37").Value = 102867000
$ws.Range("B38").Value = 12700000
$ws.Range("B39").Value = 6806648.659999877
$ws.Range("B40").Value = 250000
$ws.Range("B41").Value = 6530000
$ws.Range("B43").Value = 6780000
$ws.Range("B47").Value = 3220000

# B19 changes from a numeric value to a text value "29,966k"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "29,966k"

# Text (percentage-style, stored as literal text) value changes
$textCells = @{
    "B8"  = "22.80%"
    "B9"  = "N/A"
    "B10" = "7.84%"
    "B11" = "28.04%"
    "B12" = "6.49%"
    "B13" = "0.31%"
    "B14" = "2.33%"
    "B15" = "3.44%"
    "B16" = "2.61%"
    "B21" = "1.59%"
    "B22" = "0.51%"
    "B23" = "8.57%"
    "B26" = "0.62%"
}

foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $textCells[$addr]
}
